$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: "Tonal_Taw_HMS_LR" (tonal annoyance weighting option)
$ws.Range("F1").Value = "Tonal_Taw_HMS_LR"
$ws.Range("F2").Value = 1.03
$ws.Range("F3").Value = 1.37
$ws.Range("F4").Value = 0.154
$ws.Range("F5").Value = 0.167

# Match header formatting used by the other header cells (e.g. A1/C1/E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the selection where the author left it after adding the values
$ws.Range("F4:F5").Select() | Out-Null
